# Updates cryptos list (prices / 1h volume %) per GitHub Actions refresh.
# Note: price values that look like plain numbers (single decimal point,
# e.g. "113.90") are written with a leading apostrophe so Excel keeps them
# as text and preserves exact formatting (matching trailing zeros etc.),
# just like the original workbook stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.974.28"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "2.647.26"
$ws.Range("E3").Value = "  +6.39%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'113.90"
$ws.Range("E5").Value = "  +8.48%  "
$ws.Range("D6").Value = "'326.85"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +4.48%  "
$ws.Range("D10").Value = "'41.06"
$ws.Range("E10").Value = "  +6.41%  "
$ws.Range("D11").Value = "'20.19"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'0.0823"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'7.40"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").Value = "3.061.36"
$ws.Range("E15").Value = "  +6.25%  "
$ws.Range("D16").Value = "2.651.90"
$ws.Range("E16").Value = "  +6.37%  "
$ws.Range("E17").Value = "  +5.82%  "
$ws.Range("D18").Value = "49.883.82"
$ws.Range("D19").Value = "'13.21"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "'6.79"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").Value = "'72.08"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").Value = "'276.72"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("D26").Value = "'26.85"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +5.48%  "
$ws.Range("D31").Value = "'0.141"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").Value = "'50.30"
$ws.Range("D33").Value = "'5.46"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "'5.02"
$ws.Range("E37").Value = "  +10.09%  "
$ws.Range("E38").Value = "  +7.47%  "
$ws.Range("E39").Value = "  +8.82%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'124.03"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.113"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "'22.13"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("D45").Value = "2.084.51"
$ws.Range("E45").Value = "  +4.56%  "
$ws.Range("D46").Value = "'3.32"
$ws.Range("E46").Value = "  +7.15%  "
$ws.Range("E47").Value = "  +15.97%  "
$ws.Range("E48").Value = "  +6.24%  "
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("E50").Value = "  +5.25%  "
$ws.Range("D51").Value = "'59.76"
$ws.Range("E51").Value = "  +6.55%  "
